$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e. before
#    the existing "2022-Q3" sheet), and fill it with the latest quarter data.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $q4.Cells.Item(1, $col).Value = $headers[$i]
}

# Fund rows: code, name, scale, total stock position, position ratio, held value, rank
$rows = @(
    @("001672", "国寿安保智慧生活股票", "10.45", "90.56", "2.57", "0.2686", 10),
    @("004818", "国寿安保目标策略灵活配置混合A", "2.76", "59.92", "2.73", "0.0753", 9),
    @("004819", "国寿安保目标策略灵活配置混合C", "1.30", "59.92", "2.73", "0.0355", 9),
    @("015921", "申万菱信国证2000指数增强A", "0.21", "94.00", "0.51", "0.0011", 7),
    @("015922", "申万菱信国证2000指数增强C", "0.08", "94.00", "0.51", "0.0004", 7)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $r + 2
    $q4.Cells.Item($row, 1).Value = $r
    # Columns B-G hold text (fund code + numeric-looking strings kept as text,
    # matching the other quarter sheets), column H is a real number (rank).
    $q4.Range($q4.Cells.Item($row, 2), $q4.Cells.Item($row, 7)).NumberFormat = "@"
    $q4.Cells.Item($row, 2).Value = $rows[$r][0]
    $q4.Cells.Item($row, 3).Value = $rows[$r][1]
    $q4.Cells.Item($row, 4).Value = $rows[$r][2]
    $q4.Cells.Item($row, 5).Value = $rows[$r][3]
    $q4.Cells.Item($row, 6).Value = $rows[$r][4]
    $q4.Cells.Item($row, 7).Value = $rows[$r][5]
    $q4.Cells.Item($row, 8).Value = $rows[$r][6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row right under the
#    header for the 2022-Q4 totals, pushing all existing rows down by one.
#    Copy the whole A2:D7 block down to A3:D8 in one shot so formatting
#    (the bordered/bold style on column A) rides along with the values,
#    then overwrite row 2 with the new quarter's totals.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2:D7").Copy($total.Range("A3:D8"))

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.38

# Renumber the index column (A) for every data row (0-based, per the
# existing convention).
for ($r = 2; $r -le 8; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

